$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Danger Fever"
$ws.Range("C7").Value = "dangers"
$ws.Range("D7").Value = "selected(fever)"

$ws.Range("B8").Value = "Danger error"
$ws.Range("C8").Value = "dangers"
$ws.Range("D8").Value = "selected-at(abc)"

$ws.Range("B9").Value = "Danger error"
$ws.Range("C9").Value = "dangers"
$ws.Range("D9").Value = "selected(fever"

$ws.Range("D10").Select()
